$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 04:35"

# --- Row 55 / 56 swap: Noruega (was row55) <-> Kazajistan (was row56), with Kazajistan's data refreshed ---
$ws.Range("A55").Value = "Kazajistan"
$ws.Range("B55").Value = 8531
$ws.Range("C55").Value = 209
$ws.Range("D55").Value = 4352
$ws.Range("E55").Value = 4144
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 35

$ws.Range("A56").Value = "Noruega"
$ws.Range("B56").Value = 8352
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 7727
$ws.Range("E56").Value = 390
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 235

# --- Rows 75/76/77 rotation: Guatemala jumps to row75 (refreshed data), Guinea->row76, Uzbekistan->row77 ---
$ws.Range("A75").Value = "Guatemala"
$ws.Range("B75").Value = 3424
$ws.Range("C75").Value = 370
$ws.Range("D75").Value = 258
$ws.Range("E75").Value = 3108
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 3
$ws.Range("H75").Value = 58

$ws.Range("A76").Value = "Guinea"
$ws.Range("B76").Value = 3275
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 1673
$ws.Range("E76").Value = 1582
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 20

$ws.Range("A77").Value = "Uzbekistan"
$ws.Range("B77").Value = 3164
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 2565
$ws.Range("E77").Value = 586
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 13

# --- Rows 199/200 swap: Belice <-> Santa Lucia (tied totals, only D/H differ) ---
$ws.Range("A199").Value = "Belice"
$ws.Range("D199").Value = 16
$ws.Range("H199").Value = 2

$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

# --- Rows 207/208 swap: Groenlandia <-> Islas Turcas y Caicos (tied totals, only D/H differ) ---
$ws.Range("A207").Value = "Groenlandia"
$ws.Range("D207").Value = 11
$ws.Range("H207").Value = 0

$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 10
$ws.Range("H208").Value = 1
